# --- Mise a jour de l'application --------------------------------------
# Adds 12 new "Entrainement" rows (training session J-3, 2025-10-08) to the
# bottom of the player-tracking table on the active sheet, mirroring the
# upstream commit "Mise a jour de l'application" (data/DonneesGPSPropres.xlsx).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Type | B Date | C Periode | D MD-marker | E Nom du joueur | F Poste
#         | G Temps joue | H..V numeric GPS metrics
$newRows = @(
    @("Entrainement", 45938, "Global", "J-3", "Romain Thunet", "center back", "01:39:59", 7.81, 0.48, 7.32, 0.39, 0.09, 0.01, 0.0, 1, 4.68, 26.64, 4.72, 37, 7, 35, 10),
    @("Entrainement", 45938, "Global", "J-3", "Amine Taiar", "center back", "01:39:34", 6.95, 0.34, 6.6, 0.29, 0.04, 0.01, 0.0, 1, 4.1, 27.9, 4.56, 30, 4, 17, 2),
    @("Entrainement", 45938, "Global", "J-3", "Mattheo Haon", "right back", "01:40:08", 7.83, 0.69, 7.12, 0.58, 0.1, 0.02, 0.0, 2, 4.65, 29.1, 5.1, 42, 17, 23, 15),
    @("Entrainement", 45938, "Global", "J-3", "Ilan Ihaddadene", "center midfield", "01:40:00", 8.77, 0.71, 8.05, 0.55, 0.15, 0.02, 0.0, 3, 5.18, 26.53, 4.45, 55, 7, 29, 5),
    @("Entrainement", 45938, "Global", "J-3", "Karahali Souaré", "right forward", "01:39:59", 7.68, 0.58, 7.08, 0.43, 0.14, 0.03, 0.0, 2, 4.07, 28.9, 5.04, 56, 11, 44, 18),
    @("Entrainement", 45938, "Global", "J-3", "Omar Benyounes", "center midfield", "01:39:17", 8.37, 0.69, 7.66, 0.55, 0.16, 0.0, 0.0, 1, 4.98, 25.34, 4.24, 43, 3, 34, 9),
    @("Entrainement", 45938, "Global", "J-3", "Ilyes Boughanmi", "center forward", "01:38:01", 6.55, 0.19, 6.35, 0.15, 0.05, 0.0, 0.0, 0, 3.17, 24.75, 4.03, 20, 1, 21, 11),
    @("Entrainement", 45938, "Global", "J-3", "Levy Ndoutoume", "left back", "01:39:17", 7.53, 0.6, 6.9, 0.45, 0.13, 0.04, 0.01, 6, 3.97, 30.6, 6.48, 93, 33, 63, 32),
    @("Entrainement", 45938, "Global", "J-3", "Emmanuel Valey", "left forward", "01:30:30", 6.56, 0.6, 5.94, 0.46, 0.13, 0.03, 0.0, 3, 3.96, 29.8, 5.33, 64, 18, 40, 21),
    @("Entrainement", 45938, "Global", "J-3", "Yoan Zouma", "center back", "01:38:27", 5.79, 0.25, 5.54, 0.22, 0.03, 0.0, 0.0, 1, 3.46, 25.81, 4.04, 23, 1, 15, 2),
    @("Entrainement", 45938, "Global", "J-3", "Malik Boussaid", "right back", "01:40:16", 8.09, 0.52, 7.55, 0.46, 0.08, 0.0, 0.0, 1, 4.24, 25.58, 5.17, 49, 8, 56, 21),
    @("Entrainement", 45938, "Global", "J-3", "Sofiane Belle", "left forward", "01:38:52", 6.53, 0.49, 6.02, 0.42, 0.08, 0.0, 0.0, 1, 3.83, 25.21, 4.51, 22, 3, 23, 6)
)

$startRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1
$numCols = 22

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $arr = New-Object 'object[,]' 1,$numCols
    for ($c = 0; $c -lt $numCols; $c++) {
        $arr[0, $c] = $rowData[$c]
    }
    $targetRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, $numCols))
    $targetRange.Value = $arr
    $ws.Cells.Item($r, 2).NumberFormat = "m/d/yy"
    $ws.Cells.Item($r, 4).HorizontalAlignment = -4108
}

$lastRow = $startRow + $newRows.Length - 1
$ws.Range("E" + ($lastRow + 5)).Select() | Out-Null

